$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 362.22223
$ws.Range("I33").Value = 362.22223
$ws.Range("K33").Value = 362.22223
$ws.Range("M33").Value = -133.22223

$ws.Range("H125").Value = 100923.4
$ws.Range("I125").Value = 334009.66
$ws.Range("J125").Value = 1029.2858
$ws.Range("K125").Value = 3006086.94
$ws.Range("L125").Value = 9263.572200000001
$ws.Range("M125").Value = -3003626.94
$ws.Range("N125").Value = -14183.5722

$ws.Range("H132").Value = 3180.2173
$ws.Range("I132").Value = 3237.25
$ws.Range("K132").Value = 9711.75
$ws.Range("M132").Value = -7181.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9361.808999999999
$ws.Range("I32").Value = 9485.440000000001
$ws.Range("J32").Value = 8551.333000000001
$ws.Range("K32").Value = 9485.440000000001
$ws.Range("L32").Value = 8551.333000000001
$ws.Range("M32").Value = -9198.440000000001
$ws.Range("N32").Value = -9125.333000000001

$ws.Range("H74").Value = 1180.1613
$ws.Range("I74").Value = 1051.4
$ws.Range("J74").Value = 1716.6666
$ws.Range("K74").Value = 1051.4
$ws.Range("L74").Value = 1716.6666
$ws.Range("M74").Value = -177.4000000000001
$ws.Range("N74").Value = -3464.6666

$ws.Range("H77").Value = 1180.1613
$ws.Range("I77").Value = 1051.4
$ws.Range("J77").Value = 1716.6666
$ws.Range("K77").Value = 5257
$ws.Range("L77").Value = 8583.333000000001
$ws.Range("M77").Value = -889
$ws.Range("N77").Value = -17319.333

$ws.Range("H110").Value = 1697.1666
$ws.Range("I110").Value = 1496.6666
$ws.Range("J110").Value = 1897.6666
$ws.Range("K110").Value = 1496.6666
$ws.Range("L110").Value = 1897.6666
$ws.Range("M110").Value = 548.3334
$ws.Range("N110").Value = -5987.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9331.157999999999
$ws.Range("I4").Value = 8185.7144
$ws.Range("J4").Value = 9999.333000000001
$ws.Range("K4").Value = 8185.7144
$ws.Range("L4").Value = 9999.333000000001
$ws.Range("M4").Value = -8073.7144
$ws.Range("N4").Value = -10223.333

$ws.Range("H31").Value = 4390.909
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 4390.909
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 4390.909
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -4980.909

$ws.Range("H34").Value = 4390.909
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 4390.909
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 4390.909
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -4794.909

$ws.Range("H50").Value = 34450
$ws.Range("J50").Value = 34450
$ws.Range("L50").Value = 34450
$ws.Range("N50").Value = -35700

$ws.Range("H99").Value = 2541.5386
$ws.Range("I99").Value = 2480.625
$ws.Range("J99").Value = 2639
$ws.Range("K99").Value = 2480.625
$ws.Range("L99").Value = 2639
$ws.Range("M99").Value = -982.625
$ws.Range("N99").Value = -5635

$ws.Range("H126").Value = 2541.5386
$ws.Range("I126").Value = 2480.625
$ws.Range("J126").Value = 2639
$ws.Range("K126").Value = 7441.875
$ws.Range("L126").Value = 7917
$ws.Range("M126").Value = -4971.875
$ws.Range("N126").Value = -12857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 471
$ws.Range("I5").Value = 463.0909
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 1389.2727
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = -1277.2727
$ws.Range("N5").Value = -1724

$ws.Range("H75").Value = 5575.75
$ws.Range("I75").Value = 671.3333
$ws.Range("J75").Value = 7210.5557
$ws.Range("K75").Value = 2013.9999
$ws.Range("L75").Value = 21631.6671
$ws.Range("M75").Value = -1015.9999
$ws.Range("N75").Value = -23627.6671

$ws.Range("H78").Value = 5575.75
$ws.Range("I78").Value = 671.3333
$ws.Range("J78").Value = 7210.5557
$ws.Range("K78").Value = 6041.9997
$ws.Range("L78").Value = 64895.0013
$ws.Range("M78").Value = -1049.9997
$ws.Range("N78").Value = -74879.0013

$ws.Range("H122").Value = 801.0625
$ws.Range("I122").Value = 492.55554
$ws.Range("J122").Value = 1197.7142
$ws.Range("K122").Value = 4432.99986
$ws.Range("L122").Value = 10779.4278
$ws.Range("M122").Value = -1982.99986
$ws.Range("N122").Value = -15679.4278

$ws.Range("H135").Value = 471
$ws.Range("I135").Value = 463.0909
$ws.Range("J135").Value = 500
$ws.Range("K135").Value = 4167.8181
$ws.Range("L135").Value = 4500
$ws.Range("M135").Value = -1632.8181
$ws.Range("N135").Value = -9570

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 761.2
$ws.Range("I107").Value = 827.0625
$ws.Range("J107").Value = 497.75
$ws.Range("K107").Value = 827.0625
$ws.Range("L107").Value = 497.75
$ws.Range("M107").Value = 1092.9375
$ws.Range("N107").Value = -4337.75

$ws.Range("H109").Value = 14836.5
$ws.Range("J109").Value = 14836.5
$ws.Range("L109").Value = 14836.5
$ws.Range("N109").Value = -16916.5

$ws.Range("H132").Value = 2342.1943
$ws.Range("I132").Value = 1722.174
$ws.Range("J132").Value = 3439.1538
$ws.Range("K132").Value = 5166.522
$ws.Range("L132").Value = 10317.4614
$ws.Range("M132").Value = -2636.522
$ws.Range("N132").Value = -15377.4614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 3000
$ws.Range("I2").Value = 3000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -2888
$ws.Range("N2").ClearContents()

$ws.Range("H13").Value = 2249.8333
$ws.Range("I13").Value = 3800
$ws.Range("J13").Value = 699.6667
$ws.Range("K13").Value = 3800
$ws.Range("L13").Value = 699.6667
$ws.Range("M13").Value = -3660
$ws.Range("N13").Value = -979.6667

$ws.Range("H20").Value = 26670.166
$ws.Range("J20").Value = 15005.25
$ws.Range("L20").Value = 15005.25
$ws.Range("N20").Value = -15485.25

$ws.Range("H34").Value = 28000
$ws.Range("I34").Value = 28000
$ws.Range("K34").Value = 28000
$ws.Range("M34").Value = -27797

$ws.Range("H109").Value = 31999.6
$ws.Range("J109").Value = 31999.6
$ws.Range("L109").Value = 31999.6
$ws.Range("N109").Value = -34773.6
